$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column Q (existing last year column) into the new column R
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new 2020 column with its values
$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 5
$ws.Range("R6").Value = 3.5
$ws.Range("R7").Value = 1.8
$ws.Range("R8").Value = 24.4
$ws.Range("R9").Value = 7.2
$ws.Range("R10").Value = 2.9
$ws.Range("R11").Value = 7.4
$ws.Range("R12").Value = 4
$ws.Range("R13").Value = 3.2
$ws.Range("R14").Value = 3.5

# Update the sheet selection to match the author's last selection (new column)
$ws.Range("R4:R14").Select() | Out-Null
